$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 569, shifting all the
# existing rows (old 569..673) down to (new 571..675).
$ws.Rows.Item(569).Insert()
$ws.Rows.Item(569).Insert()

# New row 569: same shape as the surrounding records, with fresh values.
$ws.Cells.Item(569, 1).Value = 3
$ws.Cells.Item(569, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(569, 3).Value = "Coquimbo"
$ws.Cells.Item(569, 4).Value = 44694
$ws.Cells.Item(569, 5).Value = 5
$ws.Cells.Item(569, 6).Value = 100112023
$ws.Cells.Item(569, 7).Value = "Brócoli"
$ws.Cells.Item(569, 8).Value = "Sin especificar"
$ws.Cells.Item(569, 9).Value = "Primera"
$ws.Cells.Item(569, 10).Value = 2250
$ws.Cells.Item(569, 11).Value = 900
$ws.Cells.Item(569, 12).Value = 950
$ws.Cells.Item(569, 13).Value = 924
$ws.Cells.Item(569, 14).Value = "$/unidad"
$ws.Cells.Item(569, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(569, 16).Value = 924
$ws.Cells.Item(569, 17).Value = 1
$ws.Cells.Item(569, 18).Value = "Hortaliza"

# New row 570: paired "Segunda" record for the same new date.
$ws.Cells.Item(570, 1).Value = 3
$ws.Cells.Item(570, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(570, 3).Value = "Coquimbo"
$ws.Cells.Item(570, 4).Value = 44694
$ws.Cells.Item(570, 5).Value = 5
$ws.Cells.Item(570, 6).Value = 100112023
$ws.Cells.Item(570, 7).Value = "Brócoli"
$ws.Cells.Item(570, 8).Value = "Sin especificar"
$ws.Cells.Item(570, 9).Value = "Segunda"
$ws.Cells.Item(570, 10).Value = 900
$ws.Cells.Item(570, 11).Value = 700
$ws.Cells.Item(570, 12).Value = 700
$ws.Cells.Item(570, 13).Value = 700
$ws.Cells.Item(570, 14).Value = "$/unidad"
$ws.Cells.Item(570, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(570, 16).Value = 700
$ws.Cells.Item(570, 17).Value = 1
$ws.Cells.Item(570, 18).Value = "Hortaliza"
